# Added New scripts for OPQA-4225 & OPQA-4224
# Row 12 (TCID "DRA005") gets its "Jira id" (B) and "Description" (C) cells
# extended with the new OPQA-4225 test case, appended after "||".
#
# NOTE on ordering: the new shared-string table must end with the merged
# Description text BEFORE the merged Jira id text (that's the order the
# original commit's xlsx shows), so write column C before column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDescription = 'Verify that error message " Incorrect password. Please try again."should be displayed when user enters incorrect password for existing steam account.|| Verify that when user''s account is locked due to 10 invalid authentications of existing password,user becomes locked, the user is signed out'
$newJiraId = "OPQA-4221 || OPQA-4225"

$ws.Range("C12").Value = $newDescription
$ws.Range("B12").Value = $newJiraId

# The combined description now wraps across more lines, so the row grows.
$ws.Rows.Item(12).RowHeight = 60

# Move the active selection/view like the author's saved state.
$ws.Range("C8").Select()
